# Apply KHL Injuries_Master_Clubs daily refresh (2025-12-17 run)
# Sheets: snapshot (current active injuries), returned (players back today),
# new_injured (players newly added today).
$wb = $excel.ActiveWorkbook

# --- Sheet "snapshot": full refresh of the current injuries table ---
$snapshot = $wb.Worksheets.Item("snapshot")

$snapshotData = @(
    @('АВТ','Автомобилист','avtomobilist','Да Коста Стефан','77','нападающий','22668','1369_АВТ_дакостастефан','injured_active','https://www.khl.ru/clubs/avtomobilist/team/','2025-12-17T03:01:41.244683+00:00'),
    @('АВТ','Автомобилист','avtomobilist','Денежкин Максим','24','нападающий','25396','1369_АВТ_денежкинмаксим','injured_active','https://www.khl.ru/clubs/avtomobilist/team/','2025-12-17T03:01:41.244721+00:00'),
    @('АДМ','Адмирал','admiral','Дарьин Александр','41','нападающий','26282','1369_АДМ_дарьиналександр','injured_active','https://www.khl.ru/clubs/admiral/team/','2025-12-17T03:01:44.058946+00:00'),
    @('АДМ','Адмирал','admiral','Мальцев Максим В','86','защитник','39898','1369_АДМ_мальцевмаксимв','injured_active','https://www.khl.ru/clubs/admiral/team/','2025-12-17T03:01:44.058976+00:00'),
    @('АДМ','Адмирал','admiral','Педан Руслан','53','защитник','25353','1369_АДМ_педанруслан','injured_active','https://www.khl.ru/clubs/admiral/team/','2025-12-17T03:01:44.058994+00:00'),
    @('АКБ','Ак Барс','ak_bars','Яруллин Альберт','33','защитник','16365','1369_АКБ_яруллинальберт','injured_active','https://www.khl.ru/clubs/ak_bars/team/','2025-12-17T03:01:46.293409+00:00'),
    @('АМР','Амур','amur','Абросимов Роман','94','защитник','17968','1369_АМР_абросимовроман','injured_active','https://www.khl.ru/clubs/amur/team/','2025-12-17T03:01:48.630725+00:00'),
    @('БАР','Барыс','barys','Уотерспун Тайлер','26','защитник','45769','1369_БАР_уотерспунтайлер','injured_active','https://www.khl.ru/clubs/barys/team/','2025-12-17T03:01:51.413832+00:00'),
    @('БАР','Барыс','barys','Шил Адам','41','вратарь','45713','1369_БАР_шиладам','injured_active','https://www.khl.ru/clubs/barys/team/','2025-12-17T03:01:51.413862+00:00'),
    @('ДИН','Динамо М','dynamo_msk','Готовец Кирилл','41','защитник','16034','1369_ДИН_готовецкирилл','injured_active','https://www.khl.ru/clubs/dynamo_msk/team/','2025-12-17T03:01:53.717744+00:00'),
    @('ЛАД','Лада','lada','Обидин Андрей','58','нападающий','21310','1369_ЛАД_обидинандрей','injured_active','https://www.khl.ru/clubs/lada/team/','2025-12-17T03:01:58.796986+00:00'),
    @('ЛАД','Лада','lada','Ожгихин Алексей','43','нападающий','23021','1369_ЛАД_ожгихиналексей','injured_active','https://www.khl.ru/clubs/lada/team/','2025-12-17T03:01:58.797014+00:00'),
    @('ЛОК','Локомотив','lokomotiv','Сергеев Андрей','99','защитник','15416','1369_ЛОК_сергеевандрей','injured_active','https://www.khl.ru/clubs/lokomotiv/team/','2025-12-17T03:02:01.512930+00:00'),
    @('ММГ','Металлург Мг','metallurg_mg','Сиряцкий Александр','74','защитник','42458','1369_ММГ_сиряцкийалександр','injured_active','https://www.khl.ru/clubs/metallurg_mg/team/','2025-12-17T03:02:04.251469+00:00'),
    @('НХК','Нефтехимик','neftekhimik','Пастухов Илья','55','защитник','31247','1369_НХК_пастуховилья','injured_active','https://www.khl.ru/clubs/neftekhimik/team/','2025-12-17T03:02:07.125250+00:00'),
    @('СЕВ','Северсталь','severstal','Ващенко Григорий','16','защитник','14155','1369_СЕВ_ващенкогригорий','injured_active','https://www.khl.ru/clubs/severstal/team/','2025-12-17T03:02:09.804541+00:00'),
    @('СЕВ','Северсталь','severstal','Смирнов Егор Д','47','нападающий','40906','1369_СЕВ_смирновегорд','injured_active','https://www.khl.ru/clubs/severstal/team/','2025-12-17T03:02:09.804571+00:00'),
    @('СЕВ','Северсталь','severstal','Фомин Макар','77','защитник','42087','1369_СЕВ_фоминмакар','injured_active','https://www.khl.ru/clubs/severstal/team/','2025-12-17T03:02:09.804593+00:00'),
    @('СЕВ','Северсталь','severstal','Шостак Константин','1','вратарь','27876','1369_СЕВ_шостакконстантин','injured_active','https://www.khl.ru/clubs/severstal/team/','2025-12-17T03:02:09.804611+00:00'),
    @('СИБ','Сибирь','sibir','Першаков Александр','28','нападающий','42479','1369_СИБ_першаковалександр','injured_active','https://www.khl.ru/clubs/sibir/team/','2025-12-17T03:02:12.037789+00:00'),
    @('СИБ','Сибирь','sibir','Пьянов Валентин','45','нападающий','16195','1369_СИБ_пьяноввалентин','injured_active','https://www.khl.ru/clubs/sibir/team/','2025-12-17T03:02:12.037817+00:00'),
    @('СКА','СКА','ska','Зайцев Никита И','22','защитник','16024','1369_СКА_зайцевникитаи','injured_active','https://www.khl.ru/clubs/ska/team/','2025-12-17T03:02:14.343861+00:00'),
    @('СКА','СКА','ska','Короткий Матвей','71','нападающий','41566','1369_СКА_короткийматвей','injured_active','https://www.khl.ru/clubs/ska/team/','2025-12-17T03:02:14.343888+00:00'),
    @('СКА','СКА','ska','Мёрфи Тревор','8','защитник','34733','1369_СКА_мерфитревор','injured_active','https://www.khl.ru/clubs/ska/team/','2025-12-17T03:02:14.343907+00:00'),
    @('СОЧ','ХК Сочи','hc_sochi','Бикмуллин Рафаэль','24','нападающий','22424','1369_СОЧ_бикмуллинрафаэль','injured_active','https://www.khl.ru/clubs/hc_sochi/team/','2025-12-17T03:02:17.129905+00:00'),
    @('СПР','Спартак','spartak','Вишневский Дмитрий','55','защитник','15299','1369_СПР_вишневскийдмитрий','injured_active','https://www.khl.ru/clubs/spartak/team/','2025-12-17T03:02:19.463781+00:00'),
    @('СПР','Спартак','spartak','Коростелёв Никита','35','нападающий','22149','1369_СПР_коростелевникита','injured_active','https://www.khl.ru/clubs/spartak/team/','2025-12-17T03:02:19.463811+00:00'),
    @('СПР','Спартак','spartak','Порядин Павел','24','нападающий','19258','1369_СПР_порядинпавел','injured_active','https://www.khl.ru/clubs/spartak/team/','2025-12-17T03:02:19.463828+00:00'),
    @('СЮЛ','Салават Юлаев','salavat_yulaev','Алалыкин Данил','61','нападающий','34493','1369_СЮЛ_алалыкинданил','injured_active','https://www.khl.ru/clubs/salavat_yulaev/team/','2025-12-17T03:02:21.749217+00:00'),
    @('ТОР','Торпедо','torpedo','Кручинин Алексей','78','нападающий','16355','1369_ТОР_кручининалексей','injured_active','https://www.khl.ru/clubs/torpedo/team/','2025-12-17T03:02:24.546535+00:00'),
    @('ТОР','Торпедо','torpedo','Принс Шэйн','18','нападающий','19060','1369_ТОР_принсшэйн','injured_active','https://www.khl.ru/clubs/torpedo/team/','2025-12-17T03:02:24.546563+00:00'),
    @('ЦСК','ЦСКА','cska','Бучельников Дмитрий','72','нападающий','39102','1369_ЦСК_бучельниковдмитрий','injured_active','https://www.khl.ru/clubs/cska/team/','2025-12-17T03:02:29.054012+00:00'),
    @('ЦСК','ЦСКА','cska','Моисеев Данила','93','нападающий','23931','1369_ЦСК_моисеевданила','injured_active','https://www.khl.ru/clubs/cska/team/','2025-12-17T03:02:29.054040+00:00'),
    @('ШДР','Драконы','kunlun','Бишофф Джейк','28','защитник','45490','1369_ШДР_бишоффджейк','injured_active','https://www.khl.ru/clubs/kunlun/team/','2025-12-17T03:02:31.292692+00:00'),
    @('ШДР','Драконы','kunlun','Гроло Жереми','75','защитник','45343','1369_ШДР_гроложереми','injured_active','https://www.khl.ru/clubs/kunlun/team/','2025-12-17T03:02:31.292718+00:00')
)

# Old table occupied rows 2:41; clear it fully before writing the refreshed rows.
$snapshot.Range("A2:K41").ClearContents()

# Columns E (number) and G (player_id_khl) look numeric but must stay text,
# matching the source export format. Only format the rows we are about to fill.
$snapshot.Range("E2:E36").NumberFormat = "@"
$snapshot.Range("G2:G36").NumberFormat = "@"

$r = 2
foreach ($row in $snapshotData) {
    $snapshot.Cells.Item($r, 1).Value = $row[0]
    $snapshot.Cells.Item($r, 2).Value = $row[1]
    $snapshot.Cells.Item($r, 3).Value = $row[2]
    $snapshot.Cells.Item($r, 4).Value = $row[3]
    $snapshot.Cells.Item($r, 5).Value = $row[4]
    $snapshot.Cells.Item($r, 6).Value = $row[5]
    $snapshot.Cells.Item($r, 7).Value = $row[6]
    $snapshot.Cells.Item($r, 8).Value = $row[7]
    $snapshot.Cells.Item($r, 9).Value = $row[8]
    $snapshot.Cells.Item($r, 10).Value = $row[9]
    $snapshot.Cells.Item($r, 11).Value = $row[10]
    $r = $r + 1
}

# --- Sheet "returned": players who came back from injury today ---
$returned = $wb.Worksheets.Item("returned")

$returnedData = @(
    @('АВТ','Автомобилист','Зборовский Сергей','1369_АВТ_зборовскийсергей','RETURN','2025-12-17T11:02:31.805217+08:00','2025-12-17'),
    @('АДМ','Адмирал','Цыба Арсений','1369_АДМ_цыбаарсений','RETURN','2025-12-17T11:02:31.805217+08:00','2025-12-17'),
    @('НХК','Нефтехимик','Хлыстов Никита','1369_НХК_хлыстовникита','RETURN','2025-12-17T11:02:31.805217+08:00','2025-12-17'),
    @('СЕВ','Северсталь','Цицюра Владислав','1369_СЕВ_цицюравладислав','RETURN','2025-12-17T11:02:31.805217+08:00','2025-12-17'),
    @('СКА','СКА','Зыков Валентин','1369_СКА_зыковвалентин','RETURN','2025-12-17T11:02:31.805217+08:00','2025-12-17'),
    @('СКА','СКА','Сапего Сергей','1369_СКА_сапегосергей','RETURN','2025-12-17T11:02:31.805217+08:00','2025-12-17'),
    @('СОЧ','ХК Сочи','Венгрыжановский Денис','1369_СОЧ_венгрыжановскийденис','RETURN','2025-12-17T11:02:31.805217+08:00','2025-12-17'),
    @('СОЧ','ХК Сочи','Хёфенмайер Ноэль','1369_СОЧ_хефенмайерноэль','RETURN','2025-12-17T11:02:31.805217+08:00','2025-12-17')
)

# Old table occupied rows 2:2; clear it fully before writing the refreshed rows.
$returned.Range("A2:G2").ClearContents()

# Column G (changed_day, e.g. 2025-12-17) must stay text, not become a date serial.
$returned.Range("G2:G9").NumberFormat = "@"

$r = 2
foreach ($row in $returnedData) {
    $returned.Cells.Item($r, 1).Value = $row[0]
    $returned.Cells.Item($r, 2).Value = $row[1]
    $returned.Cells.Item($r, 3).Value = $row[2]
    $returned.Cells.Item($r, 4).Value = $row[3]
    $returned.Cells.Item($r, 5).Value = $row[4]
    $returned.Cells.Item($r, 6).Value = $row[5]
    $returned.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# --- Sheet "new_injured": newly injured players today ---
$newInjured = $wb.Worksheets.Item("new_injured")

$newInjuredData = @(
    @('АДМ','Адмирал','Мальцев Максим В','1369_АДМ_мальцевмаксимв','INJURED_NEW','2025-12-17T11:02:31.805217+08:00','2025-12-17'),
    @('БАР','Барыс','Шил Адам','1369_БАР_шиладам','INJURED_NEW','2025-12-17T11:02:31.805217+08:00','2025-12-17'),
    @('НХК','Нефтехимик','Пастухов Илья','1369_НХК_пастуховилья','INJURED_NEW','2025-12-17T11:02:31.805217+08:00','2025-12-17')
)

# Column G (changed_day, e.g. 2025-12-17) must stay text, not become a date serial.
$newInjured.Range("G2:G4").NumberFormat = "@"

$r = 2
foreach ($row in $newInjuredData) {
    $newInjured.Cells.Item($r, 1).Value = $row[0]
    $newInjured.Cells.Item($r, 2).Value = $row[1]
    $newInjured.Cells.Item($r, 3).Value = $row[2]
    $newInjured.Cells.Item($r, 4).Value = $row[3]
    $newInjured.Cells.Item($r, 5).Value = $row[4]
    $newInjured.Cells.Item($r, 6).Value = $row[5]
    $newInjured.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
